$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.488422666666668
$ws.Range("H2").Value = 19.465268
$ws.Range("I2").Value = 0.3444204430827928
$ws.Range("J2").Value = 0.361825925629615
$ws.Range("M2").Value = 26.07194833333334
$ws.Range("N2").Value = 78.215845
$ws.Range("O2").Value = 0.7595928012803946
$ws.Range("P2").Value = 0.7818221335728009
$ws.Range("Q2").Value = 169.1658205301623
$ws.Range("R2").Value = 1522.49238477146
$ws.Range("S2").Value = 0.2616192891794933
$ws.Range("T2").Value = 0.2828835171576992
$ws.Range("G3").Value = 6.488422666666668
$ws.Range("H3").Value = 19.465268
$ws.Range("I3").Value = 0.3444204430827928
$ws.Range("J3").Value = 0.361825925629615
$ws.Range("O3").Value = 0.155109029208254
$ws.Range("P3").Value = 0.1596482641062294
$ws.Range("Q3").Value = 34.54370045822134
$ws.Range("R3").Value = 310.893304123992
$ws.Range("S3").Value = 0.05342272056604871
$ws.Range("T3").Value = 0.05776488093539769
$ws.Range("G4").Value = 6.488422666666668
$ws.Range("H4").Value = 19.465268
$ws.Range("I4").Value = 0.3444204430827928
$ws.Range("J4").Value = 0.361825925629615
$ws.Range("M4").Value = 2.927739
$ws.Range("N4").Value = 5.855478
$ws.Range("O4").Value = 0.08529816951135136
$ws.Range("P4").Value = 0.05852960232096958
$ws.Range("Q4").Value = 18.996408089684
$ws.Range("R4").Value = 113.978448538104
$ws.Range("S4").Value = 0.02937843333725081
$ws.Range("T4").Value = 0.02117752753651808
$ws.Range("I5").Value = 0.4517209651039303
$ws.Range("J5").Value = 0.4745489404232121
$ws.Range("M5").Value = 26.07194833333334
$ws.Range("N5").Value = 78.215845
$ws.Range("O5").Value = 0.7595928012803946
$ws.Range("P5").Value = 0.7818221335728009
$ws.Range("Q5").Value = 221.8676308192139
$ws.Range("R5").Value = 1996.808677372925
$ws.Range("S5").Value = 0.3431239932803778
$ws.Range("T5").Value = 0.3710128650863877
$ws.Range("I6").Value = 0.4517209651039303
$ws.Range("J6").Value = 0.4745489404232121
$ws.Range("O6").Value = 0.155109029208254
$ws.Range("P6").Value = 0.1596482641062294
$ws.Range("S6").Value = 0.07006600037028622
$ws.Range("T6").Value = 0.07576091457201628
$ws.Range("I7").Value = 0.4517209651039303
$ws.Range("J7").Value = 0.4745489404232121
$ws.Range("M7").Value = 2.927739
$ws.Range("N7").Value = 5.855478
$ws.Range("O7").Value = 0.08529816951135136
$ws.Range("P7").Value = 0.05852960232096958
$ws.Range("Q7").Value = 24.914536776545
$ws.Range("R7").Value = 149.48722065927
$ws.Range("S7").Value = 0.03853097145326628
$ws.Range("T7").Value = 0.02777516076480809
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.331724
$ws.Range("H8").Value = 0.9951719999999999
$ws.Range("I8").Value = 0.01760867516355742
$ws.Range("J8").Value = 0.0184985395557192
$ws.Range("M8").Value = 26.07194833333334
$ws.Range("N8").Value = 78.215845
$ws.Range("O8").Value = 0.7595928012803946
$ws.Range("P8").Value = 0.7818221335728009
$ws.Range("Q8").Value = 8.648690988926667
$ws.Range("R8").Value = 77.83821890034
$ws.Range("S8").Value = 0.01337542289432309
$ws.Range("T8").Value = 0.01446256766343324
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.331724
$ws.Range("H9").Value = 0.9951719999999999
$ws.Range("I9").Value = 0.01760867516355742
$ws.Range("J9").Value = 0.0184985395557192
$ws.Range("O9").Value = 0.155109029208254
$ws.Range("P9").Value = 0.1596482641062294
$ws.Range("Q9").Value = 1.766064740152
$ws.Range("R9").Value = 15.894582661368
$ws.Range("S9").Value = 0.002731264510262885
$ws.Range("T9").Value = 0.00295325972857099
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 0.6666666666666666
$ws.Range("G10").Value = 0.331724
$ws.Range("H10").Value = 0.9951719999999999
$ws.Range("I10").Value = 0.01760867516355742
$ws.Range("J10").Value = 0.0184985395557192
$ws.Range("M10").Value = 2.927739
$ws.Range("N10").Value = 5.855478
$ws.Range("O10").Value = 0.08529816951135136
$ws.Range("P10").Value = 0.05852960232096958
$ws.Range("Q10").Value = 0.9712012920359998
$ws.Range("R10").Value = 5.827207752215999
$ws.Range("S10").Value = 0.001501987758971443
$ws.Range("T10").Value = 0.00108271216371497
$ws.Range("G11").Value = 2.7186785
$ws.Range("H11").Value = 5.437357
$ws.Range("I11").Value = 0.1443137264130649
$ws.Range("J11").Value = 0.1010711349827635
$ws.Range("M11").Value = 26.07194833333334
$ws.Range("N11").Value = 78.215845
$ws.Range("O11").Value = 0.7595928012803946
$ws.Range("P11").Value = 0.7818221335728009
$ws.Range("Q11").Value = 70.88124538694417
$ws.Range("R11").Value = 425.2874723216651
$ws.Range("S11").Value = 0.1096196677093125
$ws.Range("T11").Value = 0.07901965039484869
$ws.Range("G12").Value = 2.7186785
$ws.Range("H12").Value = 5.437357
$ws.Range("I12").Value = 0.1443137264130649
$ws.Range("J12").Value = 0.1010711349827635
$ws.Range("O12").Value = 0.155109029208254
$ws.Range("P12").Value = 0.1596482641062294
$ws.Range("Q12").Value = 14.473967028793
$ws.Range("R12").Value = 86.843802172758
$ws.Range("S12").Value = 0.02238436200535607
$ws.Range("T12").Value = 0.01613583125124458
$ws.Range("G13").Value = 2.7186785
$ws.Range("H13").Value = 5.437357
$ws.Range("I13").Value = 0.1443137264130649
$ws.Range("J13").Value = 0.1010711349827635
$ws.Range("M13").Value = 2.927739
$ws.Range("N13").Value = 5.855478
$ws.Range("O13").Value = 0.08529816951135136
$ws.Range("P13").Value = 0.05852960232096958
$ws.Range("Q13").Value = 7.9595810729115
$ws.Range("R13").Value = 31.838324291646
$ws.Range("S13").Value = 0.01230969669839639
$ws.Range("T13").Value = 0.005915653336670182
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 0.790022
$ws.Range("H14").Value = 2.370066
$ws.Range("I14").Value = 0.04193619023665445
$ws.Range("J14").Value = 0.04405545940869034
$ws.Range("M14").Value = 26.07194833333334
$ws.Range("N14").Value = 78.215845
$ws.Range("O14").Value = 0.7595928012803946
$ws.Range("P14").Value = 0.7818221335728009
$ws.Range("Q14").Value = 20.59741276619667
$ws.Range("R14").Value = 185.37671489577
$ws.Range("S14").Value = 0.03185442821688789
$ws.Range("T14").Value = 0.03444353327043221
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 0.790022
$ws.Range("H15").Value = 2.370066
$ws.Range("I15").Value = 0.04193619023665445
$ws.Range("J15").Value = 0.04405545940869034
$ws.Range("O15").Value = 0.155109029208254
$ws.Range("P15").Value = 0.1596482641062294
$ws.Range("Q15").Value = 4.205996545756
$ws.Range("R15").Value = 37.853968911804
$ws.Range("S15").Value = 0.006504681756300133
$ws.Range("T15").Value = 0.007033377618999864
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 0.790022
$ws.Range("H16").Value = 2.370066
$ws.Range("I16").Value = 0.04193619023665445
$ws.Range("J16").Value = 0.04405545940869034
$ws.Range("M16").Value = 2.927739
$ws.Range("N16").Value = 5.855478
$ws.Range("O16").Value = 0.08529816951135136
$ws.Range("P16").Value = 0.05852960232096958
$ws.Range("Q16").Value = 2.312978220258
$ws.Range("R16").Value = 13.877869321548
$ws.Range("S16").Value = 0.003577080263466429
$ws.Range("T16").Value = 0.002578548519258263
